# Update the 2025 row (row 8) of the metricas_recorrencia_anual sheet
# with refreshed BIBI recurrence metrics (commit: "atualizei dados para BIBI e ADD 06-05-2025")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 805
$ws.Range("D8").Value = 141
$ws.Range("E8").Value = 664
$ws.Range("F8").Value = 5.783429040196883
$ws.Range("G8").Value = 82.48447204968944
$ws.Range("H8").Value = 17.51552795031056
